# Adds two new worksheets - "Table S6 - HostVsymb Plast AIC" and
# "Table S7 - HostVsymb Plast GLM" - at the end of the workbook, each
# populated with the HostVsymb-by-plasticity PERMANOVA summary table
# (identical table layout/content for both sheets, mirroring the existing
# "Table S5 - HostVsymb PERMANOVA" sheet).

$wb = $excel.ActiveWorkbook

# The shared table contents (header row + 15 data rows, columns A-L).
$tableData = @(
    @(' ', 'Df', 'Sum of Squares', 'R2', 'F', 'P-value', 'Df ', 'Sum of Squares ', 'R2 ', 'F ', 'P-value ', 'species'),
    @('pCO2', 3, 1, 0.019, 0.55, 0.74217, 3, 61056, 0.208, 8.15, 0.00067, 'SSID'),
    @('temperature', 1, 3, 0.075, 6.65, 0.00333, 1, 7468, 0.025, 2.99, 0.09927, 'SSID'),
    @('reef environment', 1, 0, 0.006, 0.52, 0.56762, 1, 24705, 0.084, 9.9, 0.002, 'SSID'),
    @('Residual', 80, 30, 0.901, $null, $null, 80, 199684, 0.682, $null, $null, 'SSID'),
    @('Total', 85, 34, 1, $null, $null, 85, 292913, 1, $null, $null, 'SSID'),
    @('pCO2', 3, 1, 0.041, 1.23, 0.30513, 3, 26899, 0.024, 1.31, 0.27981, 'PSTR'),
    @('temperature', 1, 3, 0.147, 13.12, 0.00067, 1, 515173, 0.456, 75.24, 0.00067, 'PSTR'),
    @('reef environment', 1, 0, 0.02, 1.75, 0.15656, 1, 101793, 0.09, 14.87, 0.00067, 'PSTR'),
    @('Residual', 71, 14, 0.793, $null, $null, 71, 486140, 0.43, $null, $null, 'PSTR'),
    @('Total', 76, 18, 1, $null, $null, 76, 1130005, 1, $null, $null, 'PSTR'),
    @('pCO2', 3, 2, 0.136, 3.48, 0.01532, 3, 29037, 0.205, 7.04, 0.00133, 'PAST'),
    @('temperature', 1, 0, 0.036, 2.76, 0.08195, 1, 26338, 0.186, 19.15, 0.00067, 'PAST'),
    @('reef environment', 1, 0, 0.021, 1.64, 0.17988, 1, 724, 0.005, 0.53, 0.47901, 'PAST'),
    @('Residual', 62, 10, 0.807, $null, $null, 62, 85288, 0.603, $null, $null, 'PAST'),
    @('Total', 67, 13, 1, $null, $null, 67, 141387, 1, $null, $null, 'PAST')
)

# Column widths (Excel "characters" units) used on the existing, similarly
# shaped "Table S5" sheet. Sheet6 defines all 7 columns; sheet7 only the
# first 6 (matching the target workbook's per-sheet column metadata).
$colWidths6 = @(16.71, 3.71, 14.71, 5.71, 5.71, 7.71, 3.71)
$colWidths7 = @(16.71, 3.71, 14.71, 5.71, 5.71, 7.71)

function Fill-Table($ws, $colWidths) {
    for ($r = 0; $r -lt $tableData.Count; $r++) {
        $row = $tableData[$r]
        for ($c = 0; $c -lt $row.Count; $c++) {
            $val = $row[$c]
            if ($null -ne $val) {
                $ws.Cells.Item($r + 1, $c + 1).Value = $val
            }
        }
    }
    for ($c = 0; $c -lt $colWidths.Count; $c++) {
        $ws.Columns.Item($c + 1).ColumnWidth = $colWidths[$c]
    }
    $ws.PageSetup.PaperSize = 9
    $ws.PageSetup.Orientation = 1
}

# Append the two new sheets, after the current last sheet, in order.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add($null, $lastSheet)
$ws6.Name = "Table S6 - HostVsymb Plast AIC"

$ws7 = $wb.Worksheets.Add($null, $ws6)
$ws7.Name = "Table S7 - HostVsymb Plast GLM"

Fill-Table $ws6 $colWidths6
Fill-Table $ws7 $colWidths7

# Match the view zoom on the new sheets.
$ws6.Activate()
$excel.ActiveWindow.Zoom = 100
$ws7.Activate()
$excel.ActiveWindow.Zoom = 100

# Restore the originally active sheet/tab selection.
$wb.Worksheets.Item(1).Activate()

Write-Output "Added sheets: $($ws6.Name), $($ws7.Name)"
